$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Hide columns J, K, L (10, 11, 12) - "teilweise Spalten ausblenden"
# ---------------------------------------------------------------------------
$ws.Columns.Item(10).Hidden = $true
$ws.Columns.Item(11).Hidden = $true
$ws.Columns.Item(12).Hidden = $true

# ---------------------------------------------------------------------------
# 2) Re-layout the header row 7/8 for columns J:L so the single-row labels
#    "Nettoarbeitstage Monat" / "Nettoarbeitstage Intervall" / "Anteil des
#    Monats" move from row 8 up to row 7 and the cells become vertically
#    merged (J7:J8, K7:K8, L7:L8), matching the other header columns.
# ---------------------------------------------------------------------------

# Remember the label text currently sitting in row 8 before we move it.
$labelJ = $ws.Range("J8").Value2
$labelK = $ws.Range("K8").Value2
$labelL = $ws.Range("L8").Value2

# Apply the "top half of a vertical merge" look (fill + border without a
# bottom line) that is already used by Q7:S7, then recolor the alignment to
# centered to build the new style used by the merged J7/K7/L7 cells.
$ws.Range("Q7").Copy()
$ws.Range("J7:L7").PasteSpecial(-4122)
$ws.Range("J7:L7").HorizontalAlignment = -4108
$ws.Range("J7:L7").VerticalAlignment = -4160

# Apply the "bottom half of a vertical merge" look (fill + border without a
# top line), again recolored to centered, matching Q8:S8.
$ws.Range("Q8").Copy()
$ws.Range("J8:L8").PasteSpecial(-4122)
$ws.Range("J8:L8").HorizontalAlignment = -4108
$ws.Range("J8:L8").VerticalAlignment = -4160

# Move the label text up into row 7 and clear row 8.
$ws.Range("J7").Value = $labelJ
$ws.Range("K7").Value = $labelK
$ws.Range("L7").Value = $labelL
$ws.Range("J8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()

# Merge the header cells vertically, like all of their neighbouring columns.
$ws.Range("J7:J8").Merge()
$ws.Range("K7:K8").Merge()
$ws.Range("L7:L8").Merge()

# ---------------------------------------------------------------------------
# 3) Update the stored selection to match the author's final cursor position.
# ---------------------------------------------------------------------------
$ws.Range("M13").Select()
